# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4738.4614
$ws.Range("I51").Value = 2016.6666
$ws.Range("J51").Value = 7071.4287
$ws.Range("K51").Value = 2016.6666
$ws.Range("L51").Value = 7071.4287
$ws.Range("M51").Value = -1532.6666
$ws.Range("N51").Value = -8039.4287

$ws.Range("H64").Value = 3991.1333
$ws.Range("I64").Value = 3823.2
$ws.Range("J64").Value = 4327
$ws.Range("K64").Value = 3823.2
$ws.Range("L64").Value = 4327
$ws.Range("M64").Value = -3575.2
$ws.Range("N64").Value = -4823

$ws.Range("H67").Value = 3991.1333
$ws.Range("I67").Value = 3823.2
$ws.Range("J67").Value = 4327
$ws.Range("K67").Value = 3823.2
$ws.Range("L67").Value = 4327
$ws.Range("M67").Value = -2965.2
$ws.Range("N67").Value = -6043

$ws.Range("H103").Value = 361837.56
$ws.Range("I103").Value = 700.25
$ws.Range("J103").Value = 650747.4
$ws.Range("K103").Value = 2100.75
$ws.Range("L103").Value = 1952242.2
$ws.Range("M103").Value = -1514.75
$ws.Range("N103").Value = -1953414.2

$ws.Range("H111").Value = 504500.5
$ws.Range("J111").Value = 504500.5
$ws.Range("L111").Value = 1513501.5
$ws.Range("N111").Value = -1519635.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2221.4
$ws.Range("I88").Value = 2150
$ws.Range("J88").Value = 2507
$ws.Range("K88").Value = 2150
$ws.Range("L88").Value = 2507
$ws.Range("M88").Value = -1744
$ws.Range("N88").Value = -3319

$ws.Range("H91").Value = 2221.4
$ws.Range("I91").Value = 2150
$ws.Range("J91").Value = 2507
$ws.Range("K91").Value = 2150
$ws.Range("L91").Value = 2507
$ws.Range("M91").Value = -746
$ws.Range("N91").Value = -5315

$ws.Range("H97").Value = 1107.75
$ws.Range("I97").Value = 765
$ws.Range("J97").Value = 1222
$ws.Range("K97").Value = 765
$ws.Range("L97").Value = 1222
$ws.Range("M97").Value = -269
$ws.Range("N97").Value = -2214

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2229.6667
$ws.Range("I94").Value = 1828.091
$ws.Range("J94").Value = 2860.7144
$ws.Range("K94").Value = 1828.091
$ws.Range("L94").Value = 2860.7144
$ws.Range("M94").Value = -1377.091
$ws.Range("N94").Value = -3762.7144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6000
$ws.Range("J80").Value = 6000
$ws.Range("L80").Value = 18000
$ws.Range("N80").Value = -19872

$ws.Range("H83").Value = 6000
$ws.Range("J83").Value = 6000
$ws.Range("L83").Value = 54000
$ws.Range("N83").Value = -63360

$ws.Range("H124").Value = 7218.2856
$ws.Range("I124").Value = 2030
$ws.Range("J124").Value = 8083
$ws.Range("K124").Value = 6090
$ws.Range("L124").Value = 24249
$ws.Range("M124").Value = -1180
$ws.Range("N124").Value = -34069

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 35106.6
$ws.Range("J15").Value = 35106.6
$ws.Range("L15").Value = 35106.6
$ws.Range("N15").Value = -35682.6

$ws.Range("H70").Value = 5756.0654
$ws.Range("J70").Value = 5414.615
$ws.Range("L70").Value = 5414.615
$ws.Range("N70").Value = -5954.615

$ws.Range("H73").Value = 5756.0654
$ws.Range("J73").Value = 5414.615
$ws.Range("L73").Value = 5414.615
$ws.Range("N73").Value = -7286.615

$ws.Range("H81").Value = 35106.6
$ws.Range("J81").Value = 35106.6
$ws.Range("L81").Value = 35106.6
$ws.Range("N81").Value = -37102.6

$ws.Range("H84").Value = 35106.6
$ws.Range("J84").Value = 35106.6
$ws.Range("L84").Value = 105319.8
$ws.Range("N84").Value = -115303.8

$ws.Range("H132").Value = 3386.1538
$ws.Range("I132").Value = 2412.25
$ws.Range("J132").Value = 3819
$ws.Range("K132").Value = 7236.75
$ws.Range("L132").Value = 11457
$ws.Range("M132").Value = -4706.75
$ws.Range("N132").Value = -16517

$ws.Range("H138").Value = 35000
$ws.Range("J138").Value = 35000
$ws.Range("L138").Value = 35000
$ws.Range("N138").Value = -45280

$ws.Range("H140").Value = 38527.273
$ws.Range("J140").Value = 38527.273
$ws.Range("L140").Value = 38527.273
$ws.Range("N140").Value = -48887.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 166668510
$ws.Range("I68").Value = 1410
$ws.Range("J68").Value = 333335600
$ws.Range("K68").Value = 1410
$ws.Range("L68").Value = 333335600
$ws.Range("M68").Value = -661
$ws.Range("N68").Value = -333337098

$ws.Range("H71").Value = 166668510
$ws.Range("I71").Value = 1410
$ws.Range("J71").Value = 333335600
$ws.Range("K71").Value = 7050
$ws.Range("L71").Value = 1666678000
$ws.Range("M71").Value = -3306
$ws.Range("N71").Value = -1666685488

$ws.Range("H100").Value = 2769.6667
$ws.Range("J100").Value = 2820.6667
$ws.Range("L100").Value = 2820.6667
$ws.Range("N100").Value = -3902.6667

$ws.Range("H122").Value = 5091111
$ws.Range("I122").Value = 5104020
$ws.Range("J122").Value = 5000750
$ws.Range("K122").Value = 15312060
$ws.Range("L122").Value = 15002250
$ws.Range("M122").Value = -15309610
$ws.Range("N122").Value = -15007150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 34667.668
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 50001.5
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 50001.5
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -51249.5

$ws.Range("H65").Value = 34667.668
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 50001.5
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 250007.5
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -256247.5

$ws.Range("H96").Value = 4212
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 4212
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 4212
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -6958

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
